# LOQ4257.xlsx content/layout fix
# - rewrites the "Objetivos:" Portuguese body text
# - corrects the "Docentes responsáveis:" row (value had landed one row too high)
# - re-aligns the Programa resumido / Short syllabus / Programa / Syllabus block
#   (labels were offset by one row against their bodies)
# - corrects Método / Critério / Norma de recuperação / Bibliografia block
# - adds a new Bibliografia body paragraph + row 22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos: body text (PT) gets the real objectives paragraph ---
$ws.Cells.Item(10,2).Value = "Fornecer uma visão geral sobre gestão de projetos em uma organização, suas etapas e inter-relação com as demais áreasorganizacionais e sua importância no mundo competitivo dos negócios."
$ws.Cells.Item(10,3).Value = "Fornecer uma visão geral sobre gestão de projetos em uma organização, suas etapas e inter-relação com as demais áreasorganizacionais e sua importância no mundo competitivo dos negócios."

# --- Row 13: drop the "Programa resumido:" label (moves to row 14) and place
#     the professor's name (previously mis-placed on row 18) here instead ---
$ws.Cells.Item(13,1).Clear()
$ws.Cells.Item(13,2).Value = "5840917 - Fabricio Maciel Gomes"
$ws.Cells.Item(13,3).Value = "5840917 - Fabricio Maciel Gomes"
$ws.Rows.Item(13).EntireRow.AutoFit()

# --- Row 14: "Programa resumido:" + its real (Portuguese) short syllabus ---
$ws.Cells.Item(14,1).Value = "Programa resumido:"
$ws.Cells.Item(14,2).Value = "Considerações gerais sobre gerenciamento de projetos, Iniciação de projetos, Planejamento e Plano de Gerenciamento,Estrutura de Monitoramento e Avaliação, Execução e Controle."
$ws.Cells.Item(14,3).Value = "Considerações gerais sobre gerenciamento de projetos, Iniciação de projetos, Planejamento e Plano de Gerenciamento,Estrutura de Monitoramento e Avaliação, Execução e Controle."

# --- Row 15: "Short syllabus:" + its English body (was wrongly holding a date) ---
$ws.Cells.Item(15,1).Value = "Short syllabus:"
$ws.Cells.Item(15,2).Value = "General considerations about project management, Project initiation, Planning and Management Plan, Monitoring and EvaluationStructure, Execution and Control."
$ws.Cells.Item(15,3).Value = "General considerations about project management, Project initiation, Planning and Management Plan, Monitoring and EvaluationStructure, Execution and Control."
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16: "Programa:" + full Portuguese programme (same text as row 14 body) ---
$ws.Cells.Item(16,1).Value = "Programa:"
$ws.Cells.Item(16,2).Value = "Considerações gerais sobre gerenciamento de projetos, Iniciação de projetos, Planejamento e Plano de Gerenciamento,Estrutura de Monitoramento e Avaliação, Execução e Controle."
$ws.Cells.Item(16,3).Value = "Considerações gerais sobre gerenciamento de projetos, Iniciação de projetos, Planejamento e Plano de Gerenciamento,Estrutura de Monitoramento e Avaliação, Execução e Controle."

# --- Row 17: "Syllabus:" + full English programme (new B/C content on this row) ---
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Cells.Item(17,1).Value = "Syllabus:"
$ws.Cells.Item(17,2).Value = "General considerations about project management, Project initiation, Planning and Management Plan, Monitoring and EvaluationStructure, Execution and Control."
$ws.Cells.Item(17,3).Value = "General considerations about project management, Project initiation, Planning and Management Plan, Monitoring and EvaluationStructure, Execution and Control."
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18: becomes just the "Avaliação:" label, body cells cleared ---
$ws.Cells.Item(18,1).Value = "Avaliação:"
$ws.Range("B18:C18").Clear()
$ws.Rows.Item(18).EntireRow.AutoFit()

# --- Row 19: "Método:" (body text unchanged) ---
$ws.Cells.Item(19,1).Value = "Método:"

# --- Row 20: "Critério:" (body text unchanged) ---
$ws.Cells.Item(20,1).Value = "Critério:"

# --- Row 21: "Norma de recuperação:" (body text unchanged), shorter row height now ---
$ws.Cells.Item(21,1).Value = "Norma de recuperação:"
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22 (new): "Bibliografia:" + bibliography paragraph ---
$ws.Range("A21:C21").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Cells.Item(22,1).Value = "Bibliografia:"
$ws.Cells.Item(22,2).Value = "1. PMBOK. Um Guia Do Conhecimento em Gerenciamento de projetos. 5 ed. Project Management Institute. 20122. CLAUSING, D. Total quality development a step by step guide to world class concurrent engineering. New York: ASME Press,1994.3. MEREDITH, J R; MANTEL, S J; WILEY, J. Project Management: a managerial approach. 1995.4. MAXIMIANO, A . C. Administração de projetos, Atlas: São Paulo, 1997.5. SHTUB, A BARD J. F. e GLOBERSON S. Project management, Prentice hall, 1994."
$ws.Cells.Item(22,3).Value = "1. PMBOK. Um Guia Do Conhecimento em Gerenciamento de projetos. 5 ed. Project Management Institute. 20122. CLAUSING, D. Total quality development a step by step guide to world class concurrent engineering. New York: ASME Press,1994.3. MEREDITH, J R; MANTEL, S J; WILEY, J. Project Management: a managerial approach. 1995.4. MAXIMIANO, A . C. Administração de projetos, Atlas: São Paulo, 1997.5. SHTUB, A BARD J. F. e GLOBERSON S. Project management, Prentice hall, 1994."
$ws.Rows.Item(22).RowHeight = 120

# --- Column layout: column A no longer shares col-1/2 width grouping ---
$ws.Columns.Item(1).ColumnWidth = 30.7109375
